$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the two "Boton 'Regresar': ... pantalla principal." list
#    items (Historia 32 / Historia 39 task lists). A third, unrelated
#    "Regresar" bullet ("... ventana lista de productos.") must stay.
# ------------------------------------------------------------------
$targetText = [char]0x201C + "Regresar" + [char]0x201D + ": Creaci" + [char]0x00F3 + "n de un bot" + [char]0x00F3 + "n que nos redirija a la pantalla principal."

$toDelete = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$targetText*") {
        $toDelete += $p
    }
}

# Delete from the last match to the first so earlier (still-pending)
# paragraph references are not invalidated by upstream shifts.
for ($idx = $toDelete.Count - 1; $idx -ge 0; $idx--) {
    $toDelete[$idx].Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new empty (bold, 11pt) paragraph right after the table
#    for Historia 32 (the one that used to contain the first removed
#    bullet), before the pre-existing blank paragraph that follows it.
# ------------------------------------------------------------------
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

foreach ($tbl in $d.Tables) {
    $tblText = $tbl.Range.Text
    if ($tblText -like "*1 hora*" -and $tblText -like "*H32*") {
        $insertAt = $d.Range($tbl.Range.End, $tbl.Range.End)
        $insertAt.InsertXML($newParaXml)
        break
    }
}

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
